# Welcome.pptx - "fixed little issues with code / navigation formatting"
#
# The Agenda table ("Table 2") on slide 2 has a few time-range typos
# where the end time was missing its leading zero (or was simply a
# wrong hour). Fix them up so the displayed ranges read consistently.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the agenda table shape by name rather than a hard-coded index.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}

$tbl = $tableShape.Table

# Row 2, Col 1: 08h30-9h00   -> 08h30-09h00
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "08h30-09h00"

# Row 3, Col 1: 09h00-9h15   -> 09h00-09h15
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "09h00-09h15"

# Row 12, Col 1: 15h15-12h30 -> 15h15-16h30 (wrong hour typo'd as 12h30)
$tbl.Cell(12, 1).Shape.TextFrame.TextRange.Text = "15h15-16h30"
